$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E5").Value = "EUR / t"
$ws.Range("E6:E7").Value = "EUR / kg"

$ws.Range("F9").Select()
